$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.467318296432495
$ws.Range("B1").Value = 3.312880516052246
$ws.Range("C1").Value = 4.190833568572998
$ws.Range("D1").Value = 2.358716726303101
$ws.Range("E1").Value = 0.7248721122741699
